$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing value
$ws.Range("D18").Value = 0.6447426901493167

# Update existing value
$ws.Range("C19").Value = 0.2386249091493167

# Add new value
$ws.Range("D19").Value = 0.597740902

# Update existing value
$ws.Range("B20").Value = -0.0107480648506833

# Add new value
$ws.Range("C20").Value = 0.042359665
